# Add the season "record" columns (Wins / Losses / Ties) to the roster sheet.
# Previously the sheet only tracked per-player stats (A:AC); this appends three
# new columns (AD:AF) with the team's season win/loss/tie record, filled in
# on every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells should look like the existing headers (bold, bordered,
# centered) -- copy the formatting from the last existing header cell (AC1)
# onto the three new header cells before writing their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player (rows 2-51) shares the same team season record: 80-82-0.
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 80
    $ws.Cells.Item($r, 31).Value = 82
    $ws.Cells.Item($r, 32).Value = 0
}
